# Actualización desde MV -datos-
# Adds new daily "Tasa de Politica Monetaria" rows for 28-09-2021 .. 04-10-2021,
# and fills in several previously-blank cells on the last existing row (193,
# 27-09-2021) that were backfilled once more data became available.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Backfill the existing last row (193 = 27-09-2021) ----
$ws.Range("B193").Value = 0.25
$ws.Range("D193").Value = -0.1
$ws.Range("E193").Value = -0.5
$ws.Range("H193").Value = 1.75
$ws.Range("M193").Value = 1.13
$ws.Range("O193").Value = 38
$ws.Range("P193").Value = 6.25
$ws.Range("Q193").Value = 1.75
$ws.Range("R193").Value = 4.5
$ws.Range("S193").Value = 1

# ---- Row 194 = 28-09-2021 ----
$ws.Range("A194").Value = "28-09-2021"
$ws.Range("B194").Value = 0.25
$ws.Range("C194").Value = 0.1
$ws.Range("D194").Value = -0.1
$ws.Range("E194").Value = -0.5
$ws.Range("F194").Value = 0.75
$ws.Range("G194").Value = 4.35
$ws.Range("H194").Value = 1.75
$ws.Range("I194").Value = 0.1
$ws.Range("K194").Value = 6.75
$ws.Range("L194").Value = 0.5
$ws.Range("M194").Value = 1.13
$ws.Range("N194").Value = 18
$ws.Range("O194").Value = 38
$ws.Range("P194").Value = 6.25
$ws.Range("Q194").Value = 1.75
$ws.Range("R194").Value = 4.5
$ws.Range("S194").Value = 1

# ---- Row 195 = 29-09-2021 ----
$ws.Range("A195").Value = "29-09-2021"
$ws.Range("B195").Value = 0.25
$ws.Range("C195").Value = 0.1
$ws.Range("D195").Value = -0.1
$ws.Range("E195").Value = -0.5
$ws.Range("F195").Value = 0.75
$ws.Range("G195").Value = 4.35
$ws.Range("H195").Value = 1.75
$ws.Range("I195").Value = 0.1
$ws.Range("J195").Value = 0.75
$ws.Range("K195").Value = 6.75
$ws.Range("L195").Value = 0.5
$ws.Range("M195").Value = 1.13
$ws.Range("N195").Value = 18
$ws.Range("O195").Value = 38
$ws.Range("P195").Value = 6.25
$ws.Range("Q195").Value = 1.75
$ws.Range("R195").Value = 4.5
$ws.Range("S195").Value = 1

# ---- Row 196 = 30-09-2021 ----
$ws.Range("A196").Value = "30-09-2021"
$ws.Range("B196").Value = 0.25
$ws.Range("C196").Value = 0.1
$ws.Range("D196").Value = -0.1
$ws.Range("E196").Value = -0.5
$ws.Range("F196").Value = 0.75
$ws.Range("G196").Value = 4.35
$ws.Range("H196").Value = 1.75
$ws.Range("I196").Value = 0.1
$ws.Range("J196").Value = 0.75
$ws.Range("K196").Value = 6.75
$ws.Range("L196").Value = 0.5
$ws.Range("M196").Value = 1.13
$ws.Range("N196").Value = 18
$ws.Range("O196").Value = 38
$ws.Range("P196").Value = 6.25
$ws.Range("Q196").Value = 1.75
$ws.Range("R196").Value = 4.75
$ws.Range("S196").Value = 1

# ---- Row 197 = 01-10-2021 ----
# NOTE: "01-10-2021" (day <= 12) is ambiguous as dd-mm-yyyy vs mm-dd-yyyy, so
# Excel's smart entry would silently convert it to a date serial. Force it in
# as literal text (leading apostrophe = text qualifier) and then clear the
# formatting flag Excel sets for that, so the cell ends up as a plain shared
# string with no cell-level style - exactly like its neighbours.
$ws.Range("A197").Value = "'01-10-2021"
$ws.Range("A197").ClearFormats()
$ws.Range("B197").Value = 0.25
$ws.Range("C197").Value = 0.1
$ws.Range("D197").Value = -0.1
$ws.Range("E197").Value = -0.5
$ws.Range("F197").Value = 0.75
$ws.Range("H197").Value = 1.75
$ws.Range("I197").Value = 0.1
$ws.Range("J197").Value = 1.5
$ws.Range("K197").Value = 6.75
$ws.Range("L197").Value = 0.5
$ws.Range("M197").Value = 1.13
$ws.Range("N197").Value = 18
$ws.Range("O197").Value = 38
$ws.Range("P197").Value = 6.25
$ws.Range("Q197").Value = 2
$ws.Range("R197").Value = 4.75
$ws.Range("S197").Value = 1

# ---- Row 198 = 04-10-2021 ----
# Same ambiguous-date issue as row 197 above.
$ws.Range("A198").Value = "'04-10-2021"
$ws.Range("A198").ClearFormats()
$ws.Range("C198").Value = 0.1
$ws.Range("F198").Value = 0.75
$ws.Range("J198").Value = 1.5
$ws.Range("K198").Value = 6.75
$ws.Range("L198").Value = 0.5
$ws.Range("N198").Value = 18
$ws.Range("R198").Value = 4.75
$ws.Range("S198").Value = 1
